$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Kerangka Acuan ANDAL " -> "Kerangka Acuan " (two spots in the document,
#    each with slightly different trailing text). Replace the longer/more
#    specific string first so the two locations are disambiguated correctly.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Kerangka Acuan ANDAL Rencana ", $true, $false, $false, $false, $false, $true, 1, $false, "Kerangka Acuan Rencana ", 1) | Out-Null
$d.Content.Find.Execute("Kerangka Acuan ANDAL ", $true, $false, $false, $false, $false, $true, 1, $false, "Kerangka Acuan ", 1) | Out-Null

# ---------------------------------------------------------------------------
# 2) Insert 26 blank paragraphs (tabs @426, justify both, Tahoma/noProof/20)
#    right before the paragraph that carries the anchored picture, i.e.
#    immediately after the last of the pre-existing plain blank paragraphs
#    that follows the highlighted blank paragraph.
# ---------------------------------------------------------------------------
$anchor = $d.Paragraphs(71)
for ($i = 0; $i -lt 26; $i++) {
    $anchor.Range.InsertParagraphAfter()
}

# ---------------------------------------------------------------------------
# 3) Style bookkeeping: mark a few built-in styles with their normal
#    ui-priority + unhideWhenUsed flags (as Word itself does when it
#    upgrades/resaves a document created by an older Word version).
# ---------------------------------------------------------------------------
$s1 = $d.Styles("Default Paragraph Font")
$s1.Priority = 1
$s1.UnhideWhenUsed = $true

$s2 = $d.Styles("Normal Table")
$s2.Priority = 99
$s2.UnhideWhenUsed = $true

$s3 = $d.Styles("No List")
$s3.Priority = 99
$s3.UnhideWhenUsed = $true
